$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 2 and row 3 and need to be swapped
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}
